$d = $word.ActiveDocument

function Remove-ParagraphContaining($doc, [string]$needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like "*$needle*") {
            $p.Range.Delete()
            return $true
        }
    }
    return $false
}

# The Archive now just embeds the PS script directly (no separate description
# paragraphs about the archive contents / what the script does), so drop both
# explanatory paragraphs under DESCRIPTION.
Remove-ParagraphContaining $d "The Archive contains Agent Procedure Folder XML-file" | Out-Null
Remove-ParagraphContaining $d "The PowerShell script is used for detecting and uninstalling Mozilla Firefox." | Out-Null

# Step "2. Upload the PowerShell file to the Shared Files directory..." is no
# longer needed (the script now travels embedded with the procedure), so the
# install steps renumber down to: 1. Extract..., 3. Import..., 4. Execute...
Remove-ParagraphContaining $d "Upload the Power Shell file to the Shared Files directory" | Out-Null
